$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: D4 formula
$ws.Range("D4").Formula = "=IF(A2=1,B4+4,C4)"

# Row 6: A6, B6, C6 values and D6/E6 formulas
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = "A"
$ws.Range("C6").Value = "B"
$ws.Range("D6").Formula = "=IF(A6=1,B6,C6)"
$ws.Range("E6").Formula = "=IF(D6=""A"",D4,0)"

# Row 4: E4 literal string "T" (added after the new strings "A"/"B" so shared-string order matches)
$ws.Range("E4").Value = "T"

# Update the selected cell as in the saved view state
$ws.Range("D5").Select()

$wb.Application.Calculate()
